$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 357.70834
$ws.Range("I33").Value = 200.94444
$ws.Range("J33").Value = 828
$ws.Range("K33").Value = 200.94444
$ws.Range("L33").Value = 828
$ws.Range("M33").Value = 28.05556000000001
$ws.Range("N33").Value = -1286
$ws.Range("H40").Value = 2600
$ws.Range("J40").Value = 2600
$ws.Range("L40").Value = 2600
$ws.Range("N40").Value = -2950
$ws.Range("H86").Value = 18232958
$ws.Range("I86").Value = 3220.7144
$ws.Range("K86").Value = 3220.7144
$ws.Range("M86").Value = -2097.7144
$ws.Range("H89").Value = 18232958
$ws.Range("I89").Value = 3220.7144
$ws.Range("K89").Value = 16103.572
$ws.Range("M89").Value = -10487.572
$ws.Range("H98").Value = 1785.8096
$ws.Range("I98").Value = 1865.1
$ws.Range("K98").Value = 1865.1
$ws.Range("M98").Value = -367.0999999999999
$ws.Range("H122").Value = 1785.8096
$ws.Range("I122").Value = 1865.1
$ws.Range("K122").Value = 5595.299999999999
$ws.Range("M122").Value = -3145.299999999999
$ws.Range("H135").Value = 47621290
$ws.Range("I135").Value = 58825850
$ws.Range("K135").Value = 529432650
$ws.Range("M135").Value = -529430115

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1122.6364
$ws.Range("I2").Value = 1167
$ws.Range("J2").Value = 971.8
$ws.Range("K2").Value = 1167
$ws.Range("L2").Value = 971.8
$ws.Range("M2").Value = -1054
$ws.Range("N2").Value = -1197.8
$ws.Range("H4").Value = 13648.889
$ws.Range("J4").Value = 17514.285
$ws.Range("L4").Value = 17514.285
$ws.Range("N4").Value = -17746.285
$ws.Range("H5").Value = 5877.143
$ws.Range("J5").Value = 10247.5
$ws.Range("L5").Value = 10247.5
$ws.Range("N5").Value = -10471.5
$ws.Range("H32").Value = 3316.05
$ws.Range("J32").Value = 3997.8333
$ws.Range("L32").Value = 3997.8333
$ws.Range("N32").Value = -4571.8333
$ws.Range("H45").Value = 1715.9231
$ws.Range("I45").Value = 1608.9166
$ws.Range("K45").Value = 1608.9166
$ws.Range("M45").Value = -1231.9166
$ws.Range("H116").Value = 1122.6364
$ws.Range("I116").Value = 1167
$ws.Range("J116").Value = 971.8
$ws.Range("K116").Value = 1167
$ws.Range("L116").Value = 971.8
$ws.Range("M116").Value = 1127
$ws.Range("N116").Value = -5559.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1122.6364
$ws.Range("I3").Value = 1167
$ws.Range("J3").Value = 971.8
$ws.Range("K3").Value = 1167
$ws.Range("L3").Value = 971.8
$ws.Range("M3").Value = -1053
$ws.Range("N3").Value = -1199.8
$ws.Range("H4").Value = 5877.143
$ws.Range("J4").Value = 10247.5
$ws.Range("L4").Value = 10247.5
$ws.Range("N4").Value = -10477.5
$ws.Range("H42").Value = 500000
$ws.Range("J42").Value = 500000
$ws.Range("L42").Value = 500000
$ws.Range("N42").Value = -500656

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5411.6665
$ws.Range("I62").Value = 3117.5
$ws.Range("K62").Value = 3117.5
$ws.Range("M62").Value = -2493.5
$ws.Range("H65").Value = 5411.6665
$ws.Range("I65").Value = 3117.5
$ws.Range("K65").Value = 15587.5
$ws.Range("M65").Value = -12467.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1973.5
$ws.Range("J5").Value = 5000
$ws.Range("L5").Value = 15000
$ws.Range("N5").Value = -15224
$ws.Range("H8").Value = 1130.0667
$ws.Range("I8").Value = 1130.0667
$ws.Range("K8").Value = 3390.2001
$ws.Range("M8").Value = -3251.2001
$ws.Range("H55").Value = 800.875
$ws.Range("J55").Value = 2747.5
$ws.Range("L55").Value = 8242.5
$ws.Range("N55").Value = -8596.5
$ws.Range("H81").Value = 7934.3706
$ws.Range("I81").Value = 932
$ws.Range("K81").Value = 2796
$ws.Range("M81").Value = -1673
$ws.Range("H84").Value = 7934.3706
$ws.Range("I84").Value = 932
$ws.Range("K84").Value = 8388
$ws.Range("M84").Value = -2772
$ws.Range("H122").Value = 767.7143
$ws.Range("I122").Value = 685.5
$ws.Range("J122").Value = 877.3333
$ws.Range("K122").Value = 6169.5
$ws.Range("L122").Value = 7895.9997
$ws.Range("M122").Value = -3719.5
$ws.Range("N122").Value = -12795.9997
$ws.Range("H135").Value = 1973.5
$ws.Range("J135").Value = 5000
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -50070
$ws.Range("H137").Value = 1549.75
$ws.Range("I137").Value = 999.5
$ws.Range("J137").Value = 2100
$ws.Range("K137").Value = 2998.5
$ws.Range("L137").Value = 6300
$ws.Range("M137").Value = 2101.5
$ws.Range("N137").Value = -16500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 56000
$ws.Range("J15").Value = 56000
$ws.Range("L15").Value = 56000
$ws.Range("N15").Value = -56576
$ws.Range("H43").Value = 145549.19
$ws.Range("I43").Value = 178737.23
$ws.Range("K43").Value = 178737.23
$ws.Range("M43").Value = -178586.23
$ws.Range("H80").Value = 2609.7058
$ws.Range("I80").Value = 2384.2222
$ws.Range("K80").Value = 2384.2222
$ws.Range("M80").Value = -1386.2222
$ws.Range("H81").Value = 56000
$ws.Range("J81").Value = 56000
$ws.Range("L81").Value = 56000
$ws.Range("N81").Value = -57996
$ws.Range("H83").Value = 2609.7058
$ws.Range("I83").Value = 2384.2222
$ws.Range("K83").Value = 11921.111
$ws.Range("M83").Value = -6929.111000000001
$ws.Range("H84").Value = 56000
$ws.Range("J84").Value = 56000
$ws.Range("L84").Value = 168000
$ws.Range("N84").Value = -177984
$ws.Range("H132").Value = 4061.111
$ws.Range("I132").Value = 3315.4614
$ws.Range("J132").Value = 5999.8
$ws.Range("K132").Value = 9946.3842
$ws.Range("L132").Value = 17999.4
$ws.Range("M132").Value = -7416.3842
$ws.Range("N132").Value = -23059.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 83335320
$ws.Range("I7").Value = 100001380
$ws.Range("K7").Value = 100001380
$ws.Range("M7").Value = -100001268
$ws.Range("H22").Value = 1045.6
$ws.Range("I22").Value = 1056.875
$ws.Range("K22").Value = 1056.875
$ws.Range("M22").Value = -761.875
$ws.Range("H27").Value = 1045.6
$ws.Range("I27").Value = 1056.875
$ws.Range("K27").Value = 1056.875
$ws.Range("M27").Value = -949.875
$ws.Range("H46").Value = 3495.5715
$ws.Range("J46").Value = 4280
$ws.Range("L46").Value = 4280
$ws.Range("N46").Value = -4656
$ws.Range("H61").Value = 36963.133
$ws.Range("I61").Value = 36034.785
$ws.Range("K61").Value = 36034.785
$ws.Range("M61").Value = -35832.785
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H113").Value = 36963.133
$ws.Range("I113").Value = 36034.785
$ws.Range("K113").Value = 36034.785
$ws.Range("M113").Value = -33864.785
$ws.Range("H126").Value = 83335320
$ws.Range("I126").Value = 100001380
$ws.Range("K126").Value = 300004140
$ws.Range("M126").Value = -300001670
$ws.Range("H136").Value = 2251.6875
$ws.Range("I136").Value = 2079.3845
$ws.Range("J136").Value = 2998.3333
$ws.Range("K136").Value = 6238.1535
$ws.Range("L136").Value = 8994.999899999999
$ws.Range("M136").Value = -3688.1535
$ws.Range("N136").Value = -14094.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 45297.75
$ws.Range("J46").Value = 45297.75
$ws.Range("L46").Value = 45297.75
$ws.Range("N46").Value = -45759.75
$ws.Range("H81").Value = 8339211.5
$ws.Range("I81").Value = 3412.3572
$ws.Range("K81").Value = 6824.7144
$ws.Range("M81").Value = -5763.7144
$ws.Range("H84").Value = 8339211.5
$ws.Range("I84").Value = 3412.3572
$ws.Range("K84").Value = 34123.572
$ws.Range("M84").Value = -28819.572
$ws.Range("H126").Value = 2083.3333
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H132").Value = 3158.611
$ws.Range("I132").Value = 3055.3547
$ws.Range("J132").Value = 3798.8
$ws.Range("K132").Value = 9166.0641
$ws.Range("L132").Value = 11396.4
$ws.Range("M132").Value = -6636.0641
$ws.Range("N132").Value = -16456.4
$ws.Range("H134").Value = 45297.75
$ws.Range("J134").Value = 45297.75
$ws.Range("L134").Value = 135893.25
$ws.Range("N134").Value = -140963.25
$ws.Range("H136").Value = 4879.5625
$ws.Range("I136").Value = 2988.625
$ws.Range("J136").Value = 6770.5
$ws.Range("K136").Value = 8965.875
$ws.Range("L136").Value = 20311.5
$ws.Range("M136").Value = -6415.875
$ws.Range("N136").Value = -25411.5
